# Apply cryptocurrency price/volume updates for Tue Feb 27 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.967.58"
$ws.Range("E2").Value = "  +8.96%  "

$ws.Range("D3").Value = "3.217.57"
$ws.Range("E3").Value = "  +4.05%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'395.98"
$ws.Range("E5").Value = "  +2.44%  "

$ws.Range("D6").Value = "'110.75"
$ws.Range("E6").Value = "  +7.07%  "

$ws.Range("E7").Value = "  +2.56%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.617"
$ws.Range("E9").Value = "  +5.50%  "

$ws.Range("D10").Value = "'39.14"
$ws.Range("E10").Value = "  +6.12%  "

$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  +6.41%  "

$ws.Range("E12").Value = "  +2.07%  "

$ws.Range("D13").Value = "3.725.03"
$ws.Range("E13").Value = "  +3.83%  "

$ws.Range("D14").Value = "'8.06"
$ws.Range("E14").Value = "  +3.80%  "

$ws.Range("D15").Value = "'19.00"
$ws.Range("E15").Value = "  +2.84%  "

$ws.Range("D16").Value = "3.220.55"
$ws.Range("E16").Value = "  +3.91%  "

$ws.Range("D17").Value = "'1.05"
$ws.Range("E17").Value = "  +5.45%  "

$ws.Range("D18").Value = "'10.89"
$ws.Range("E18").Value = "  +2.41%  "

$ws.Range("D19").Value = "55.801.71"
$ws.Range("E19").Value = "  +8.42%  "

$ws.Range("D20").Value = "'3.32"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("E21").Value = "  +5.85%  "

$ws.Range("D22").Value = "'12.94"
$ws.Range("E22").Value = "  +3.93%  "

$ws.Range("D23").Value = "'298.19"
$ws.Range("E23").Value = "  +12.24%  "

$ws.Range("D24").Value = "'75.38"
$ws.Range("E24").Value = "  +7.62%  "

$ws.Range("D25").Value = "'3.20"
$ws.Range("E25").Value = "  +1.74%  "

$ws.Range("D26").Value = "'8.14"
$ws.Range("E26").Value = "  +2.09%  "

$ws.Range("D27").Value = "'28.07"
$ws.Range("E27").Value = "  +2.96%  "

$ws.Range("D28").Value = "'7.49"
$ws.Range("E28").Value = "  +4.84%  "

$ws.Range("D29").Value = "'0.172"
$ws.Range("E29").Value = "  +4.13%  "

$ws.Range("E30").Value = "  +0.53%  "

$ws.Range("E31").Value = "  +3.69%  "

$ws.Range("D32").Value = "'11.10"
$ws.Range("E32").Value = "  +6.97%  "

$ws.Range("D33").Value = "'0.0490"
$ws.Range("E33").Value = "  +4.00%  "

$ws.Range("D34").Value = "'36.06"
$ws.Range("E34").Value = "  +1.51%  "

$ws.Range("E35").Value = "  +2.58%  "

$ws.Range("D36").Value = "'51.39"
$ws.Range("E36").Value = "  +2.86%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'3.12"
$ws.Range("E37").Value = "  +25.39%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  +4.74%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("D40").Value = "'135.65"
$ws.Range("E40").Value = "  +4.83%  "

$ws.Range("E41").Value = "  +4.76%  "

$ws.Range("E42").Value = "  +2.94%  "

$ws.Range("D43").Value = "'3.98"
$ws.Range("E43").Value = "  +4.62%  "

$ws.Range("E44").Value = "  +2.98%  "

$ws.Range("D45").Value = "'0.284"
$ws.Range("E45").Value = "  -1.86%  "

$ws.Range("D46").Value = "'22.10"
$ws.Range("E46").Value = "  +0.39%  "

$ws.Range("D47").Value = "'2.16"
$ws.Range("E47").Value = "  +52.20%  "

$ws.Range("E48").Value = "  +1.61%  "

$ws.Range("D49").Value = "'2.47"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("D50").Value = "2.126.20"
$ws.Range("E50").Value = "  +2.59%  "

$ws.Range("E51").Value = "  +7.59%  "
